# Add two new Amazon order rows at the top of the data (below the header row),
# pushing the existing rows down by two, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right after the header row (row 1).
$ws.Rows("2:3").Insert() | Out-Null

# --- Row 2: order 026-1450774-3129108 ---
$ws.Range("A2").Value = "026-1450774-3129108"
$ws.Range("C2").Value = "2025-09-11T18:20:36+00:00"
$ws.Range("D2").Value = "2025-09-11T18:58:47+00:00"
$ws.Range("E2").Value = "Pending"
$ws.Range("F2").Value = "Merchant"
$ws.Range("G2").Value = "Amazon.co.uk"
$ws.Range("H2").Value = "WebsiteOrderChannel"
$ws.Range("J2").Value = "Standard"
$ws.Range("K2").Value = "Buckchi Motion Sensor Lights Indoor, 4 Pack Cupboard Light, Stair Lights Night Light, 3 Modes USB Charging Wall, Magnetic for Kitchen Stair Closet Under Cabinet, 50 Lumens 6000K Warm White"
$ws.Range("L2").Value = "8Led Usb Light Warm 250601000010"
$ws.Range("M2").Value = "B0FND8L7N5"
$ws.Range("N2").Value = "Unshipped"
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = "GBP"
$ws.Range("Q2").Value = 14.49
$ws.Range("R2").Value = 2.42
$ws.Range("Y2").Value = "Langho"
$ws.Range("Z2").Value = "Blackburn"
$ws.Range("AA2").Value = "BB6 8DT"
$ws.Range("AB2").Value = "GB"
$ws.Range("AD2").Value = $false
$ws.Range("AG2").Value = $false

# --- Row 3: order 205-8164344-9869957 ---
$ws.Range("A3").Value = "205-8164344-9869957"
$ws.Range("C3").Value = "2025-09-11T16:24:57+00:00"
$ws.Range("D3").Value = "2025-09-11T16:55:03+00:00"
$ws.Range("E3").Value = "Pending"
$ws.Range("F3").Value = "Merchant"
$ws.Range("G3").Value = "Amazon.co.uk"
$ws.Range("H3").Value = "WebsiteOrderChannel"
$ws.Range("J3").Value = "Standard"
$ws.Range("K3").Value = "Cordless Drill Driver 21V, Wuppertal Cordless Hammer Drill Set, 45N.m Battery Drill, 25+3 Torque Set, 2-Speed, LED Light, 1500mAH Battery, 30PCS Electric Drill Set for Home and Garden DIY Project"
$ws.Range("L3").Value = "Electric Drill 2506010000001"
$ws.Range("M3").Value = "B0FLQKQLX5"
$ws.Range("N3").Value = "Unshipped"
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = "GBP"
$ws.Range("Q3").Value = 23.99
$ws.Range("R3").Value = 4
$ws.Range("Y3").Value = "Liverpool"
$ws.Range("Z3").Value = "Merseyside"
$ws.Range("AA3").Value = "L21 8HU"
$ws.Range("AB3").Value = "GB"
$ws.Range("AD3").Value = $false
$ws.Range("AG3").Value = $false

# Match the workbook's saved selection/view state after the edit.
$ws.Range("E4").Select() | Out-Null
